# The deck's single Design ("Integral" / "Red Violet") had its theme colors
# swapped with the otherwise-unused default "Office Theme" that only the
# Notes Master pointed at. Net effect on the theme actually driving every
# slide/layout/master in the deck: its 12 theme colors go back to the
# stock Office palette (the color *names* shown in the Design gallery are
# cosmetic metadata the host doesn't round-trip, so we focus on the part
# that is visibly in effect presentation-wide).

$p = $ppt.ActivePresentation

function Set-ThemeRGB($themeColorScheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColorScheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Order exposed by ThemeColorScheme.Item(n): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    Set-ThemeRGB $tcs $i $officeThemeColors[$i - 1]
}
